$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tableau1")

# --- Fix existing row 47 (Heure début / Heure fin were corrected) ---
$ws.Range("C47").Value = 0.5625
$ws.Range("D47").Value = 0.58333333333333337

# --- Grow the table by two rows (B4:L47 -> B4:L49) ---
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# Copy number formats down from row 47 so the new date/time cells reuse the
# same styles instead of creating new ones.
$ws.Range("B47").Copy()
$ws.Range("B48:B49").PasteSpecial(-4122)
$ws.Range("C47:D47").Copy()
$ws.Range("C48:D49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New unique strings must be introduced in the same order as the original
# author typed them, so they land at the same sharedStrings indices:
#   86 "Score", 87 "J'ai changer de place ..." (row 49), 88 "J'ai crée une
#   fonction ..." (row 48).
$ws.Range("H48").Value = "Score"
$ws.Range("J49").Value = "J'ai changer de place quellque variable "
$ws.Range("J48").Value = "J'ai crée une fonction qui calcule les scores et j'ai ajouter dans les regle comment le score est calculer"

# --- Row 48 : new "Score" entry ---
$ws.Range("B48").Value = 44272
$ws.Range("C48").Value = 0.58333333333333337
$ws.Range("D48").Value = 0.60416666666666663
$ws.Range("E48").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])"
$ws.Range("F48").Value = "Ma-20"
$ws.Range("G48").Value = "Code"
$ws.Range("I48").Value = "CPNV"
$ws.Range("K48").Value = "Oui"
$ws.Rows.Item(48).RowHeight = 43.2

# --- Row 49 : new "Mise au propre" entry ---
$ws.Range("B49").Value = 44272
$ws.Range("C49").Value = 0.60416666666666663
$ws.Range("D49").Value = 0.625
$ws.Range("E49").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])"
$ws.Range("F49").Value = "Ma-20"
$ws.Range("G49").Value = "Code"
$ws.Range("H49").Value = "Mise au propre"
$ws.Range("I49").Value = "CPNV"
$ws.Range("K49").Value = "oui"
$ws.Rows.Item(49).RowHeight = 28.8

# --- Match the author's final selection ---
$ws.Range("G53").Select() | Out-Null
